$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the trailing rows that no longer exist in the updated order file ---
# Previously the run had 19 trials (rows 2-20); the refreshed order only has
# 16 trials (rows 2-17), so drop the old rows 18-20 (former trials 17-19).
$ws.Range("A18:D20").EntireRow.Delete() | Out-Null

# --- New "ITI" column header ---
$ws.Range("D1").Value = "ITI"

# --- Updated ConditionType (col C) values + new ITI (col D) values, per trial ---
$updates = @(
    @{ Row = 2;  ConditionType = 1; ITI = 7 }   # Trial 1
    @{ Row = 3;  ConditionType = 1; ITI = 6 }   # Trial 2
    @{ Row = 4;  ConditionType = 3; ITI = 6 }   # Trial 3
    @{ Row = 5;  ConditionType = 4; ITI = 9 }   # Trial 4
    @{ Row = 6;  ConditionType = 2; ITI = 7 }   # Trial 5
    @{ Row = 7;  ConditionType = 2; ITI = 8 }   # Trial 6
    @{ Row = 8;  ConditionType = 3; ITI = 7 }   # Trial 7
    @{ Row = 9;  ConditionType = 2; ITI = 9 }   # Trial 8
    @{ Row = 10; ConditionType = 3; ITI = 8 }   # Trial 9
    @{ Row = 11; ConditionType = 1; ITI = 6 }   # Trial 10
    @{ Row = 12; ConditionType = 2; ITI = 6 }   # Trial 11
    @{ Row = 13; ConditionType = 4; ITI = 7 }   # Trial 12
    @{ Row = 14; ConditionType = 4; ITI = 7 }   # Trial 13
    @{ Row = 15; ConditionType = 1; ITI = 7 }   # Trial 14
    @{ Row = 16; ConditionType = 4; ITI = 6 }   # Trial 15
    @{ Row = 17; ConditionType = 3; ITI = 6 }   # Trial 16
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.ConditionType
    $ws.Cells.Item($u.Row, 4).Value = $u.ITI
}

# --- Column width for the (now wider) ConditionType column ---
$ws.Columns.Item(3).ColumnWidth = 18

# --- Selection, matching the saved view state ---
$ws.Range("I16").Select() | Out-Null
